# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The previous account-statement periods (2304..2404, rows 16-28) are
# replaced with the new ones: the period list is now sorted in descending
# order (2404..2304) instead of ascending, with each period keeping its own
# "Valor Mora" / "Salario Basico" figures. In addition, the "Salario Basico"
# database value for period 2410 (row 29) is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$firstRow = 16
$lastRow  = 28

# Snapshot the current Periodo Mora / Valor Mora / Salario Basico table
# (columns E:G) for the rows holding the monthly periods 2304..2404.
$colE = @{}
$colF = @{}
$colG = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $colE[$r] = $ws.Range("E$r").Value()
    $colF[$r] = $ws.Range("F$r").Value()
    $colG[$r] = $ws.Range("G$r").Value()
}

# Reverse the row order (2304..2404 -> 2404..2304) so every period keeps
# travelling together with its own Valor Mora / Salario Basico values.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $firstRow + $lastRow - $r
    $ws.Range("E$r").Value = $colE[$srcRow]
    $ws.Range("F$r").Value = $colF[$srcRow]
    $ws.Range("G$r").Value = $colG[$srcRow]
}

# Database update: Salario Basico for period 2410 (row 29)
$ws.Range("G29").Value = 1000000
